# Insert a new weekly price record at row 66 for
# "Hortaliza, Agrícola del Norte S.A. de Arica - Cebollín baby".
# This pushes the previous rows 66-77 down to 67-78 and fills the
# newly inserted row 66 with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(66).Insert()

$ws.Cells.Item(66, 1).Value = 1
$ws.Cells.Item(66, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(66, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(66, 4).Value = 44637
$ws.Cells.Item(66, 5).Value = 15
$ws.Cells.Item(66, 6).Value = 100112038
$ws.Cells.Item(66, 7).Value = "Cebollín baby"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 300
$ws.Cells.Item(66, 11).Value = 1400
$ws.Cells.Item(66, 12).Value = 1500
$ws.Cells.Item(66, 13).Value = 1450
$ws.Cells.Item(66, 14).Value = "`$/paquete"
$ws.Cells.Item(66, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(66, 16).Value = 1450
$ws.Cells.Item(66, 17).Value = 1
$ws.Cells.Item(66, 18).Value = "Hortaliza"
